# chore: update Sheets via scheduled runner
#
# Refresh cached market-price derived columns (currentAveragePrice /
# currentAveragePriceNQ / currentAveragePriceHQ / LevePriceNQ / LevePriceHQ /
# LeveProfitNQ / LeveProfitHQ) across the eight crafting-job sheets to match
# the latest scrape. Values below are the new cached numbers for each
# touched cell; most rows keep column types stable, with one row (WVR!139)
# flipping from an HQ-profit cell to an NQ-profit cell as the underlying
# NQ/HQ average prices swapped which side was cheaper.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 955.7
$ws.Range("J32").Value = 928.6667
$ws.Range("L32").Value = 928.6667
$ws.Range("N32").Value = -1580.6667
$ws.Range("H62").Value = 1000000000
$ws.Range("I62").Value = 1000000000
$ws.Range("K62").Value = 1000000000
$ws.Range("M62").Value = -999999376
$ws.Range("H65").Value = 1000000000
$ws.Range("I65").Value = 1000000000
$ws.Range("K65").Value = 5000000000
$ws.Range("M65").Value = -4999996880
$ws.Range("H106").Value = 2466.55
$ws.Range("I106").Value = 1888.7333
$ws.Range("K106").Value = 1888.7333
$ws.Range("M106").Value = -1257.7333
$ws.Range("H112").Value = 2669.923
$ws.Range("J112").Value = 2768.76
$ws.Range("L112").Value = 8306.280000000001
$ws.Range("N112").Value = -10522.28
$ws.Range("H138").Value = 3609.1018
$ws.Range("I138").Value = 1919.6086
$ws.Range("J138").Value = 4688.5
$ws.Range("K138").Value = 5758.825800000001
$ws.Range("L138").Value = 14065.5
$ws.Range("M138").Value = -618.8258000000005
$ws.Range("N138").Value = -24345.5
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3639468.2
$ws.Range("I32").Value = 4258158.5
$ws.Range("K32").Value = 4258158.5
$ws.Range("M32").Value = -4257871.5
$ws.Range("H97").Value = 1242846.8
$ws.Range("J97").Value = 30282.857
$ws.Range("L97").Value = 30282.857
$ws.Range("N97").Value = -31274.857
$ws.Range("H102").Value = 4023.4546
$ws.Range("I102").Value = 4023.4546
$ws.Range("K102").Value = 4023.4546
$ws.Range("M102").Value = -2401.4546
$ws.Range("H110").Value = 20402294
$ws.Range("J110").Value = 4636.273
$ws.Range("L110").Value = 4636.273
$ws.Range("N110").Value = -8726.273000000001
$ws.Range("H132").Value = 7615.357
$ws.Range("I132").Value = 4767.4736
$ws.Range("K132").Value = 14302.4208
$ws.Range("M132").Value = -11772.4208
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 595.4375
$ws.Range("I94").Value = 505
$ws.Range("K94").Value = 505
$ws.Range("M94").Value = -54
$ws.Range("H99").Value = 995
$ws.Range("I99").Value = 991.6667
$ws.Range("K99").Value = 991.6667
$ws.Range("M99").Value = 506.3333
$ws.Range("H105").Value = 47634510
$ws.Range("I105").Value = 52648372
$ws.Range("K105").Value = 52648372
$ws.Range("M105").Value = -52646625
$ws.Range("H107").Value = 2909.742
$ws.Range("I107").Value = 1669.7142
$ws.Range("K107").Value = 1669.7142
$ws.Range("M107").Value = 250.2858000000001
$ws.Range("H122").Value = 65333
$ws.Range("J122").Value = 65333
$ws.Range("L122").Value = 65333
$ws.Range("N122").Value = -75133
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H54").Value = 36748.25
$ws.Range("I54").Value = 29993
$ws.Range("K54").Value = 29993
$ws.Range("M54").Value = -29335
$ws.Range("H86").Value = 6634
$ws.Range("I86").Value = 4955.5557
$ws.Range("K86").Value = 4955.5557
$ws.Range("M86").Value = -3832.5557
$ws.Range("H89").Value = 6634
$ws.Range("I89").Value = 4955.5557
$ws.Range("K89").Value = 24777.7785
$ws.Range("M89").Value = -19161.7785
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H41").Value = 200
$ws.Range("I41").Value = 200
$ws.Range("K41").Value = 600
$ws.Range("M41").Value = -262
$ws.Range("H64").Value = 526.6667
$ws.Range("I64").Value = 540.25
$ws.Range("J64").Value = 499.5
$ws.Range("K64").Value = 1620.75
$ws.Range("L64").Value = 1498.5
$ws.Range("M64").Value = -1350.75
$ws.Range("N64").Value = -2038.5
$ws.Range("H67").Value = 526.6667
$ws.Range("I67").Value = 540.25
$ws.Range("J67").Value = 499.5
$ws.Range("K67").Value = 1620.75
$ws.Range("L67").Value = 1498.5
$ws.Range("M67").Value = -684.75
$ws.Range("N67").Value = -3370.5
$ws.Range("H68").Value = 2777.842
$ws.Range("J68").Value = 3257.276
$ws.Range("L68").Value = 9771.828
$ws.Range("N68").Value = -11393.828
$ws.Range("H71").Value = 2777.842
$ws.Range("J71").Value = 3257.276
$ws.Range("L71").Value = 29315.484
$ws.Range("N71").Value = -37427.484
$ws.Range("H107").Value = 5828.3335
$ws.Range("I107").Value = 1242.5
$ws.Range("K107").Value = 3727.5
$ws.Range("M107").Value = -1807.5
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2606.2666
$ws.Range("I68").Value = 2661.875
$ws.Range("J68").Value = 2542.7144
$ws.Range("K68").Value = 2661.875
$ws.Range("L68").Value = 2542.7144
$ws.Range("M68").Value = -1912.875
$ws.Range("N68").Value = -4040.7144
$ws.Range("H71").Value = 2606.2666
$ws.Range("I71").Value = 2661.875
$ws.Range("J71").Value = 2542.7144
$ws.Range("K71").Value = 13309.375
$ws.Range("L71").Value = 12713.572
$ws.Range("M71").Value = -9565.375
$ws.Range("N71").Value = -20201.572
$ws.Range("H100").Value = 8932320
$ws.Range("I100").Value = 13892012
$ws.Range("K100").Value = 13892012
$ws.Range("M100").Value = -13891471
$ws.Range("H122").Value = 29416456
$ws.Range("I122").Value = 45458004
$ws.Range("J122").Value = 6948.8335
$ws.Range("K122").Value = 136374012
$ws.Range("L122").Value = 20846.5005
$ws.Range("M122").Value = -136371562
$ws.Range("N122").Value = -25746.5005
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H25").Value = 64999.668
$ws.Range("J25").Value = 64999.668
$ws.Range("L25").Value = 64999.668
$ws.Range("N25").Value = -65585.66800000001
$ws.Range("H61").Value = 5046.647
$ws.Range("I61").Value = 5087.125
$ws.Range("J61").Value = 4399
$ws.Range("K61").Value = 5087.125
$ws.Range("L61").Value = 4399
$ws.Range("M61").Value = -4795.125
$ws.Range("N61").Value = -4983
$ws.Range("H122").Value = 4016.1292
$ws.Range("I122").Value = 3305.3
$ws.Range("J122").Value = 5308.5454
$ws.Range("K122").Value = 9915.900000000001
$ws.Range("L122").Value = 15925.6362
$ws.Range("M122").Value = -7465.900000000001
$ws.Range("N122").Value = -20825.6362
$ws.Range("H139").Value = 100000
$ws.Range("I139").Value = 100000
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 100000
$ws.Range("L139").Value = 0
$ws.Range("M139").Value = -94860
$ws.Range("N139").ClearContents()
